# "Flag more RO names"
#
# Insert 4 new region rows (all flagged with Issues = "RO", i.e. newly
# identified Regional-Officer name-change candidates) into the practice
# detags sheet. The rows are inserted above specific existing regions so
# the final sheet dimension grows from A1:I24 to A1:I28.
#
# Insertion points below are given as the FINAL row number each new row
# ends up at once all four inserts are applied. Performing the inserts in
# ascending order of that final row number means "insert at row N now"
# always lands the new row at the right spot, because every row that is
# still above N at the time of the insert is already in its final place.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-RegionRow($Row, $Region, $Issues, $Minor, $MinorTimestamp, $Major, $MajorTimestamp, $NativeEmbassies, $Link, $Organizations) {
    # Shift the target row (and everything below it) down by one,
    # leaving a blank row at $Row.
    $ws.Rows($Row).Insert()

    # Match the formatting used by the rest of the "Region" column
    # (bold, bordered, centered/top-aligned) instead of whatever the
    # freshly inserted blank row defaulted to.
    $ws.Cells.Item($Row, 1).Font.Bold = $true
    $ws.Cells.Item($Row, 1).HorizontalAlignment = -4108
    $ws.Cells.Item($Row, 1).VerticalAlignment = -4160
    $ws.Cells.Item($Row, 1).Borders.LineStyle = 1

    $ws.Cells.Item($Row, 1).Value = $Region
    $ws.Cells.Item($Row, 2).Value = $Issues
    $ws.Cells.Item($Row, 3).Value = $Minor
    $ws.Cells.Item($Row, 4).Value = $MinorTimestamp
    $ws.Cells.Item($Row, 5).Value = $Major
    $ws.Cells.Item($Row, 6).Value = $MajorTimestamp
    $ws.Cells.Item($Row, 7).Value = $NativeEmbassies
    $ws.Cells.Item($Row, 8).Value = $Link
    $ws.Cells.Item($Row, 9).Value = $Organizations
}

# 1) New row above "The Koprulu Sector"
Insert-RegionRow 5 "abolished" "RO" 137 "0:02:17" 206 "0:03:26" $false "https://www.nationstates.net/region=abolished" "Unknown"

# 2) New row above "The Great States of Fascist Nudists"
Insert-RegionRow 11 "Zolochiv" "RO" 421 "0:07:01" 631 "0:10:31" $false "https://www.nationstates.net/region=zolochiv" "Unknown"

# 3) New row above "Aerope"
Insert-RegionRow 13 "Propounded Empathy" "RO" 733 "0:12:13" 1099 "0:18:19" $false "https://www.nationstates.net/region=propounded_empathy" "Unknown"

# 4) Another new row, also above "Aerope" (directly below the previous one)
Insert-RegionRow 14 "Crazed Nations CN" "RO" 1087 "0:18:07" 1631 "0:27:11" $false "https://www.nationstates.net/region=crazed_nations_cn" "Unknown"
